{"js": "// Replace each arithmetic expression in the table with its updated value.\n// The document contains a single table of \"a op b=\" cells; every cell's\n// text is unique, so a literal (non-wildcard) search-and-replace for each\n// old/new pair unambiguously targets the correct cell.\nconst replacements = [\n  [\"56-25=\", \"1+97=\"],\n  [\"72-66=\", \"75-67=\"],\n  [\"33+5=\", \"29+7=\"],\n  [\"5-2=\", \"62-41=\"],\n  [\"86+11=\", \"35+60=\"],\n  [\"74-1=\", \"80-73=\"],\n  [\"8+38=\", \"57-8=\"],\n  [\"94-28=\", \"82+10=\"],\n  [\"45+3=\", \"49+1=\"],\n  [\"59+17=\", \"31+45=\"],\n  [\"87-67=\", \"20+27=\"],\n  [\"8+70=\", \"85-53=\"],\n  [\"78+9=\", \"48-42=\"],\n  [\"22+74=\", \"66-9=\"],\n  [\"8+20=\", \"22+12=\"],\n  [\"13+45=\", \"85-43=\"],\n  [\"65+7=\", \"85-75=\"],\n  [\"85-3=\", \"53-17=\"],\n  [\"34+43=\", \"6+2=\"],\n  [\"68-29=\", \"38+47=\"],\n  [\"8+10=\", \"44-4=\"],\n  [\"33-26=\", \"55-51=\"],\n  [\"68-23=\", \"27-20=\"],\n  [\"38+37=\", \"93-60=\"],\n  [\"62+23=\", \"86-48=\"],\n  [\"12+3=\", \"41-2=\"],\n  [\"31+58=\", \"51+48=\"],\n  [\"99-77=\", \"2+3=\"],\n  [\"73-33=\", \"82+1=\"],\n  [\"99-65=\", \"84-10=\"],\n  [\"92-31=\", \"53+39=\"],\n  [\"55-16=\", \"79-74=\"],\n  [\"21-15=\", \"32-23=\"],\n  [\"6+41=\", \"83-24=\"],\n  [\"68-45=\", \"25-7=\"],\n  [\"42-17=\", \"13+78=\"],\n  [\"95-61=\", \"48-20=\"],\n  [\"41+47=\", \"24+50=\"],\n  [\"72+16=\", \"66-46=\"],\n  [\"15+46=\", \"56+37=\"],\n  [\"2+57=\", \"0+27=\"],\n  [\"70+24=\", \"92-22=\"],\n  [\"54-1=\", \"41+39=\"],\n  [\"17-16=\", \"42-29=\"],\n  [\"30+26=\", \"84-22=\"],\n  [\"29-12=\", \"22+18=\"],\n  [\"92-48=\", \"9+54=\"],\n  [\"24-8=\", \"45-21=\"],\n  [\"78-77=\", \"68+31=\"],\n  [\"13+21=\", \"67-67=\"],\n  [\"41-38=\", \"6+21=\"],\n  [\"95-54=\", \"58-11=\"],\n  [\"25-11=\", \"32-31=\"],\n  [\"44-19=\", \"61-24=\"],\n  [\"48+18=\", \"8+32=\"],\n  [\"23+23=\", \"14+65=\"],\n  [\"73-19=\", \"66-48=\"],\n  [\"77-16=\", \"48-11=\"],\n  [\"3+50=\", \"82-32=\"],\n  [\"72+1=\", \"41-11=\"],\n  [\"20+38=\", \"8+48=\"],\n  [\"26+53=\", \"94-31=\"],\n  [\"5+84=\", \"24+53=\"],\n  [\"85-84=\", \"26-12=\"],\n  [\"17+52=\", \"93-61=\"],\n  [\"11+14=\", \"12+66=\"],\n  [\"16+62=\", \"80-67=\"],\n  [\"81+10=\", \"6+10=\"],\n  [\"1+84=\", \"72-25=\"],\n  [\"68-4=\", \"20+21=\"],\n  [\"60-42=\", \"12+56=\"],\n  [\"11+39=\", \"61+0=\"],\n  [\"94-67=\", \"39+5=\"],\n  [\"14+21=\", \"48+36=\"],\n  [\"65-30=\", \"34+29=\"],\n  [\"42-34=\", \"2+8=\"],\n  [\"53+40=\", \"18+78=\"],\n  [\"75-65=\", \"34+2=\"],\n  [\"13+1=\", \"5+42=\"],\n  [\"61-59=\", \"47+19=\"],\n  [\"17-1=\", \"1+30=\"],\n  [\"37+6=\", \"53-20=\"],\n  [\"27+40=\", \"7+90=\"],\n  [\"70+12=\", \"70-17=\"],\n  [\"46-21=\", \"76+15=\"],\n  [\"39-13=\", \"94-1=\"],\n  [\"58+28=\", \"30+33=\"],\n  [\"86-45=\", \"73-2=\"],\n  [\"18+53=\", \"33+43=\"],\n  [\"95+0=\", \"9+1=\"],\n  [\"55+43=\", \"46+42=\"],\n  [\"47+35=\", \"99-35=\"],\n  [\"76-19=\", \"16+83=\"],\n  [\"1+40=\", \"31+23=\"],\n  [\"35-28=\", \"4+88=\"],\n  [\"18+5=\", \"90-17=\"],\n  [\"29+67=\", \"50-23=\"],\n  [\"76-76=\", \"13+66=\"],\n  [\"33+55=\", \"24+15=\"],\n  [\"48-9=\", \"3+21=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each arithmetic expression in the table with its updated value.\n# The document contains a single table of \"a op b=\" cells; every cell's\n# text is unique, so a literal (non-wildcard) find-and-replace for each\n# old/new pair unambiguously targets the correct cell.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @{ old = '56-25='; new = '1+97=' },\n  @{ old = '72-66='; new = '75-67=' },\n  @{ old = '33+5='; new = '29+7=' },\n  @{ old = '5-2='; new = '62-41=' },\n  @{ old = '86+11='; new = '35+60=' },\n  @{ old = '74-1='; new = '80-73=' },\n  @{ old = '8+38='; new = '57-8=' },\n  @{ old = '94-28='; new = '82+10=' },\n  @{ old = '45+3='; new = '49+1=' },\n  @{ old = '59+17='; new = '31+45=' },\n  @{ old = '87-67='; new = '20+27=' },\n  @{ old = '8+70='; new = '85-53=' },\n  @{ old = '78+9='; new = '48-42=' },\n  @{ old = '22+74='; new = '66-9=' },\n  @{ old = '8+20='; new = '22+12=' },\n  @{ old = '13+45='; new = '85-43=' },\n  @{ old = '65+7='; new = '85-75=' },\n  @{ old = '85-3='; new = '53-17=' },\n  @{ old = '34+43='; new = '6+2=' },\n  @{ old = '68-29='; new = '38+47=' },\n  @{ old = '8+10='; new = '44-4=' },\n  @{ old = '33-26='; new = '55-51=' },\n  @{ old = '68-23='; new = '27-20=' },\n  @{ old = '38+37='; new = '93-60=' },\n  @{ old = '62+23='; new = '86-48=' },\n  @{ old = '12+3='; new = '41-2=' },\n  @{ old = '31+58='; new = '51+48=' },\n  @{ old = '99-77='; new = '2+3=' },\n  @{ old = '73-33='; new = '82+1=' },\n  @{ old = '99-65='; new = '84-10=' },\n  @{ old = '92-31='; new = '53+39=' },\n  @{ old = '55-16='; new = '79-74=' },\n  @{ old = '21-15='; new = '32-23=' },\n  @{ old = '6+41='; new = '83-24=' },\n  @{ old = '68-45='; new = '25-7=' },\n  @{ old = '42-17='; new = '13+78=' },\n  @{ old = '95-61='; new = '48-20=' },\n  @{ old = '41+47='; new = '24+50=' },\n  @{ old = '72+16='; new = '66-46=' },\n  @{ old = '15+46='; new = '56+37=' },\n  @{ old = '2+57='; new = '0+27=' },\n  @{ old = '70+24='; new = '92-22=' },\n  @{ old = '54-1='; new = '41+39=' },\n  @{ old = '17-16='; new = '42-29=' },\n  @{ old = '30+26='; new = '84-22=' },\n  @{ old = '29-12='; new = '22+18=' },\n  @{ old = '92-48='; new = '9+54=' },\n  @{ old = '24-8='; new = '45-21=' },\n  @{ old = '78-77='; new = '68+31=' },\n  @{ old = '13+21='; new = '67-67=' },\n  @{ old = '41-38='; new = '6+21=' },\n  @{ old = '95-54='; new = '58-11=' },\n  @{ old = '25-11='; new = '32-31=' },\n  @{ old = '44-19='; new = '61-24=' },\n  @{ old = '48+18='; new = '8+32=' },\n  @{ old = '23+23='; new = '14+65=' },\n  @{ old = '73-19='; new = '66-48=' },\n  @{ old = '77-16='; new = '48-11=' },\n  @{ old = '3+50='; new = '82-32=' },\n  @{ old = '72+1='; new = '41-11=' },\n  @{ old = '20+38='; new = '8+48=' },\n  @{ old = '26+53='; new = '94-31=' },\n  @{ old = '5+84='; new = '24+53=' },\n  @{ old = '85-84='; new = '26-12=' },\n  @{ old = '17+52='; new = '93-61=' },\n  @{ old = '11+14='; new = '12+66=' },\n  @{ old = '16+62='; new = '80-67=' },\n  @{ old = '81+10='; new = '6+10=' },\n  @{ old = '1+84='; new = '72-25=' },\n  @{ old = '68-4='; new = '20+21=' },\n  @{ old = '60-42='; new = '12+56=' },\n  @{ old = '11+39='; new = '61+0=' },\n  @{ old = '94-67='; new = '39+5=' },\n  @{ old = '14+21='; new = '48+36=' },\n  @{ old = '65-30='; new = '34+29=' },\n  @{ old = '42-34='; new = '2+8=' },\n  @{ old = '53+40='; new = '18+78=' },\n  @{ old = '75-65='; new = '34+2=' },\n  @{ old = '13+1='; new = '5+42=' },\n  @{ old = '61-59='; new = '47+19=' },\n  @{ old = '17-1='; new = '1+30=' },\n  @{ old = '37+6='; new = '53-20=' },\n  @{ old = '27+40='; new = '7+90=' },\n  @{ old = '70+12='; new = '70-17=' },\n  @{ old = '46-21='; new = '76+15=' },\n  @{ old = '39-13='; new = '94-1=' },\n  @{ old = '58+28='; new = '30+33=' },\n  @{ old = '86-45='; new = '73-2=' },\n  @{ old = '18+53='; new = '33+43=' },\n  @{ old = '95+0='; new = '9+1=' },\n  @{ old = '55+43='; new = '46+42=' },\n  @{ old = '47+35='; new = '99-35=' },\n  @{ old = '76-19='; new = '16+83=' },\n  @{ old = '1+40='; new = '31+23=' },\n  @{ old = '35-28='; new = '4+88=' },\n  @{ old = '18+5='; new = '90-17=' },\n  @{ old = '29+67='; new = '50-23=' },\n  @{ old = '76-76='; new = '13+66=' },\n  @{ old = '33+55='; new = '24+15=' },\n  @{ old = '48-9='; new = '3+21=' },\n)\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nforeach ($r in $replacements) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $r.old\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $r.new\n  $find.Execute($r.old, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $r.new, $wdReplaceAll) | Out-Null\n}\n"}
